{"js": "const NEW_VALUES = [\"25+72=\", \"77-30=\", \"98-7=\", \"72-51=\", \"20+60=\", \"53+46=\", \"39-27=\", \"26+23=\", \"32+24=\", \"28-9=\", \"47+24=\", \"15+64=\", \"56-27=\", \"36+38=\", \"80-57=\", \"60+2=\", \"1+19=\", \"36-36=\", \"7+12=\", \"28+1=\", \"15-9=\", \"81-6=\", \"19+12=\", \"74-59=\", \"97-58=\", \"44+54=\", \"66+16=\", \"1+87=\", \"96-7=\", \"24-13=\", \"66-16=\", \"8+56=\", \"20+53=\", \"10+29=\", \"82-8=\", \"46+45=\", \"42+6=\", \"98-88=\", \"29+29=\", \"28-12=\", \"9+85=\", \"60+16=\", \"39-2=\", \"89-19=\", \"1+7=\", \"10+63=\", \"82-48=\", \"86-7=\", \"88-56=\", \"8+63=\", \"35+21=\", \"4+25=\", \"43-34=\", \"18+77=\", \"49+14=\", \"79-50=\", \"15+80=\", \"8+54=\", \"26+42=\", \"34-28=\", \"3+29=\", \"65-29=\", \"30+28=\", \"69-14=\", \"13+36=\", \"88-56=\", \"84-28=\", \"61-24=\", \"50+21=\", \"34+21=\", \"59-18=\", \"86-52=\", \"80-20=\", \"3-2=\", \"81-35=\", \"45+38=\", \"5+28=\", \"87+12=\", \"58+33=\", \"37-4=\", \"41+55=\", \"32-32=\", \"19+64=\", \"43-2=\", \"50+32=\", \"4+22=\", \"86+4=\", \"88-53=\", \"71-5=\", \"34+54=\", \"14+65=\", \"38-2=\", \"79-22=\", \"40+47=\", \"82-4=\", \"67-32=\", \"31-13=\", \"27+31=\", \"52-45=\", \"26+56=\"];\nconst OLD_VALUES = [\"25+5=\", \"22+20=\", \"99-73=\", \"58+35=\", \"64-29=\", \"87-35=\", \"58-10=\", \"77-43=\", \"58+6=\", \"68-17=\", \"83-77=\", \"29-26=\", \"5+1=\", \"11+7=\", \"1+51=\", \"63-2=\", \"0+11=\", \"64+29=\", \"83-6=\", \"76-56=\", \"2+60=\", \"3+9=\", \"2+20=\", \"9+42=\", \"33+8=\", \"93-42=\", \"82-3=\", \"77-39=\", \"90-30=\", \"95+4=\", \"29+46=\", \"16+82=\", \"26-24=\", \"33-19=\", \"10+31=\", \"37+9=\", \"12+33=\", \"36-21=\", \"51+1=\", \"81-67=\", \"60+22=\", \"41-18=\", \"73+24=\", \"10+69=\", \"94-46=\", \"85-36=\", \"74-43=\", \"57-18=\", \"86-4=\", \"56-20=\", \"38-19=\", \"99-22=\", \"95-30=\", \"60-8=\", \"44-36=\", \"6+87=\", \"59+40=\", \"11+54=\", \"59+6=\", \"52-25=\", \"19+57=\", \"87-69=\", \"85-58=\", \"92-7=\", \"41-17=\", \"32+18=\", \"79-11=\", \"78-10=\", \"99-70=\", \"45-33=\", \"49+10=\", \"31-1=\", \"33+29=\", \"52-30=\", \"51+32=\", \"91-91=\", \"58-14=\", \"68+28=\", \"39-37=\", \"35+22=\", \"72-4=\", \"52+10=\", \"7+68=\", \"82-37=\", \"5+49=\", \"87-59=\", \"82+17=\", \"13+69=\", \"14+85=\", \"45-44=\", \"26+41=\", \"67+31=\", \"14+23=\", \"40+48=\", \"36+58=\", \"82-38=\", \"27+33=\", \"23-22=\", \"87-69=\", \"83-23=\"];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body.\");\n}\nconst table = tables.items[0];\ntable.load(\"rowCount, values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = table.values.length > 0 ? table.values[0].length : 0;\n\nif (rowCount * colCount !== NEW_VALUES.length) {\n  throw new Error(\n    \"Table shape \" + rowCount + \"x\" + colCount +\n    \" does not match expected \" + NEW_VALUES.length + \" cells.\"\n  );\n}\n\n// Sanity-check the current contents against what the diff expects to find,\n// in row-major order, before mutating anything.\nlet k = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const current = table.values[r][c];\n    if (current !== OLD_VALUES[k]) {\n      throw new Error(\n        \"Cell (\" + r + \",\" + c + \") text [\" + current +\n        \"] does not match expected original [\" + OLD_VALUES[k] + \"]\"\n      );\n    }\n    k++;\n  }\n}\n\n// Grab every cell's first paragraph up front (one sync) so we can\n// replace text while preserving the existing run/paragraph formatting\n// (font, size, alignment) instead of wiping it with a fresh default run.\nconst paragraphs = [];\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.body.paragraphs.load(\"items\");\n    paragraphs.push(cell.body.paragraphs);\n  }\n}\nawait context.sync();\n\nlet i = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const para = paragraphs[i].items[0];\n    const range = para.getRange();\n    range.insertText(NEW_VALUES[i], Word.InsertLocation.replace);\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "$newValues = @(\"25+72=\", \"77-30=\", \"98-7=\", \"72-51=\", \"20+60=\", \"53+46=\", \"39-27=\", \"26+23=\", \"32+24=\", \"28-9=\", \"47+24=\", \"15+64=\", \"56-27=\", \"36+38=\", \"80-57=\", \"60+2=\", \"1+19=\", \"36-36=\", \"7+12=\", \"28+1=\", \"15-9=\", \"81-6=\", \"19+12=\", \"74-59=\", \"97-58=\", \"44+54=\", \"66+16=\", \"1+87=\", \"96-7=\", \"24-13=\", \"66-16=\", \"8+56=\", \"20+53=\", \"10+29=\", \"82-8=\", \"46+45=\", \"42+6=\", \"98-88=\", \"29+29=\", \"28-12=\", \"9+85=\", \"60+16=\", \"39-2=\", \"89-19=\", \"1+7=\", \"10+63=\", \"82-48=\", \"86-7=\", \"88-56=\", \"8+63=\", \"35+21=\", \"4+25=\", \"43-34=\", \"18+77=\", \"49+14=\", \"79-50=\", \"15+80=\", \"8+54=\", \"26+42=\", \"34-28=\", \"3+29=\", \"65-29=\", \"30+28=\", \"69-14=\", \"13+36=\", \"88-56=\", \"84-28=\", \"61-24=\", \"50+21=\", \"34+21=\", \"59-18=\", \"86-52=\", \"80-20=\", \"3-2=\", \"81-35=\", \"45+38=\", \"5+28=\", \"87+12=\", \"58+33=\", \"37-4=\", \"41+55=\", \"32-32=\", \"19+64=\", \"43-2=\", \"50+32=\", \"4+22=\", \"86+4=\", \"88-53=\", \"71-5=\", \"34+54=\", \"14+65=\", \"38-2=\", \"79-22=\", \"40+47=\", \"82-4=\", \"67-32=\", \"31-13=\", \"27+31=\", \"52-45=\", \"26+56=\")\n$oldValues = @(\"25+5=\", \"22+20=\", \"99-73=\", \"58+35=\", \"64-29=\", \"87-35=\", \"58-10=\", \"77-43=\", \"58+6=\", \"68-17=\", \"83-77=\", \"29-26=\", \"5+1=\", \"11+7=\", \"1+51=\", \"63-2=\", \"0+11=\", \"64+29=\", \"83-6=\", \"76-56=\", \"2+60=\", \"3+9=\", \"2+20=\", \"9+42=\", \"33+8=\", \"93-42=\", \"82-3=\", \"77-39=\", \"90-30=\", \"95+4=\", \"29+46=\", \"16+82=\", \"26-24=\", \"33-19=\", \"10+31=\", \"37+9=\", \"12+33=\", \"36-21=\", \"51+1=\", \"81-67=\", \"60+22=\", \"41-18=\", \"73+24=\", \"10+69=\", \"94-46=\", \"85-36=\", \"74-43=\", \"57-18=\", \"86-4=\", \"56-20=\", \"38-19=\", \"99-22=\", \"95-30=\", \"60-8=\", \"44-36=\", \"6+87=\", \"59+40=\", \"11+54=\", \"59+6=\", \"52-25=\", \"19+57=\", \"87-69=\", \"85-58=\", \"92-7=\", \"41-17=\", \"32+18=\", \"79-11=\", \"78-10=\", \"99-70=\", \"45-33=\", \"49+10=\", \"31-1=\", \"33+29=\", \"52-30=\", \"51+32=\", \"91-91=\", \"58-14=\", \"68+28=\", \"39-37=\", \"35+22=\", \"72-4=\", \"52+10=\", \"7+68=\", \"82-37=\", \"5+49=\", \"87-59=\", \"82+17=\", \"13+69=\", \"14+85=\", \"45-44=\", \"26+41=\", \"67+31=\", \"14+23=\", \"40+48=\", \"36+58=\", \"82-38=\", \"27+33=\", \"23-22=\", \"87-69=\", \"83-23=\")\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\nif ($rowCount * $colCount -ne $newValues.Length) {\n    throw \"Table shape $rowCount x $colCount does not match expected $($newValues.Length) cells.\"\n}\n\n# Sanity-check the current contents (row-major order) against what the diff\n# expects to find before mutating anything. Cell.Range.Text carries a\n# trailing cell-mark (CR + BEL), so trim that off before comparing.\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($current -ne $oldValues[$i]) {\n            throw \"Cell ($r,$c) text [$current] does not match expected original [$($oldValues[$i])]\"\n        }\n        $i++\n    }\n}\n\n# Assigning Range.Text replaces the cell's text while keeping the existing\n# paragraph/run formatting (font, size, alignment) and the cell mark intact.\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
